# repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -8
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -1
$ws.Range("F17").Value = -1
$ws.Range("F20").Value = -4
$ws.Range("F28").Value = -1
